$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 400
$ws.Range("F5").Value = 5007
$ws.Range("F6").Value = 5007
$ws.Range("F7").Value = 62
$ws.Range("F9").Value = 494
$ws.Range("F11").Value = 686
$ws.Range("F12").Value = 4798
$ws.Range("F13").Value = 19
$ws.Range("F14").Value = 37
$ws.Range("F15").Value = 63
$ws.Range("F16").Value = 197
$ws.Range("F17").Value = 204
$ws.Range("F18").Value = 90
$ws.Range("F19").Value = 236
$ws.Range("F20").Value = 3704
$ws.Range("F24").Value = 3513
$ws.Range("F28").Value = 180
$ws.Range("F31").Value = 102
$ws.Range("F35").Value = 133
$ws.Range("F36").Value = 6175
$ws.Range("F37").Value = 974
$ws.Range("F38").Value = 469
$ws.Range("F40").Value = 967
$ws.Range("F42").Value = 1277
$ws.Range("F43").Value = 144
$ws.Range("F44").Value = 607
$ws.Range("F46").Value = 2139
$ws.Range("F49").Value = 746
$ws.Range("F50").Value = 890

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 17

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F5").Value = 400
$ws.Range("F6").Value = 5007
$ws.Range("F7").Value = 5007
$ws.Range("F8").Value = 62
$ws.Range("F12").Value = 494
$ws.Range("F14").Value = 686
$ws.Range("F15").Value = 4798
$ws.Range("F16").Value = 19
$ws.Range("F17").Value = 37
$ws.Range("F18").Value = 63
$ws.Range("F19").Value = 197
$ws.Range("F20").Value = 204
$ws.Range("F21").Value = 90
$ws.Range("F22").Value = 236
$ws.Range("F23").Value = 3704
$ws.Range("F24").Value = 3513
$ws.Range("F27").Value = 180
$ws.Range("F30").Value = 102
$ws.Range("F35").Value = 6175
$ws.Range("F36").Value = 974
$ws.Range("F37").Value = 469
$ws.Range("F41").Value = 967
$ws.Range("F42").Value = 1277
$ws.Range("F43").Value = 144
$ws.Range("F44").Value = 607
$ws.Range("F45").Value = 2139
$ws.Range("F48").Value = 746
$ws.Range("F49").Value = 890
